# Insert a new weekly price record at row 72 ("Hortaliza, Terminal La
# Palmera de La Serena - Albahaca"), pushing the existing rows 72:141
# down to 73:142 (dimension grows from A1:R141 to A1:R142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 72 downward (through the end of the used range) by one row.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new record.
$ws.Cells.Item(72, 1).Value = 8
$ws.Cells.Item(72, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(72, 3).Value = "Coquimbo"
$ws.Cells.Item(72, 4).Value = 44904
$ws.Cells.Item(72, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(72, 5).Value = 4
$ws.Cells.Item(72, 6).Value = 100112052
$ws.Cells.Item(72, 7).Value = "Albahaca"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 600
$ws.Cells.Item(72, 11).Value = 4000
$ws.Cells.Item(72, 12).Value = 4500
$ws.Cells.Item(72, 13).Value = 4250
$ws.Cells.Item(72, 14).Value = "$/paquete"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 4250
$ws.Cells.Item(72, 17).Value = 1
$ws.Cells.Item(72, 18).Value = "Hortaliza"
